# Auto-generated edit script: Add data for 2023-09-19
# Updates violent crime statistics (mostly 2023 column J, with a few 2021 column H corrections)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J2").Value = 5496
$ws.Range("H3").Value = 8348
$ws.Range("J3").Value = 5849
$ws.Range("J4").Value = 1272
$ws.Range("J5").Value = 449
$ws.Range("J6").Value = 7319
$ws.Range("H7").Value = 26013
$ws.Range("J7").Value = 20385

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("J2").Value = 354
$ws.Range("J6").Value = 430
$ws.Range("J7").Value = 1279

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("J6").Value = 113
$ws.Range("J7").Value = 418

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("J2").Value = 231
$ws.Range("J6").Value = 322
$ws.Range("J7").Value = 937

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("J6").Value = 71
$ws.Range("J7").Value = 305

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("J4").Value = 20
$ws.Range("J5").Value = 27
$ws.Range("J6").Value = 182
$ws.Range("J7").Value = 629

$ws = $wb.Worksheets.Item("New City")
$ws.Range("J2").Value = 154
$ws.Range("J3").Value = 150
$ws.Range("J6").Value = 185
$ws.Range("J7").Value = 525

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("J3").Value = 125
$ws.Range("J7").Value = 319

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J6").Value = 147
$ws.Range("J7").Value = 589
$ws.Range("J8").Value = 1279
$ws.Range("J11").Value = 318
$ws.Range("J13").Value = 25
$ws.Range("J14").Value = 103
$ws.Range("J18").Value = 173
$ws.Range("J19").Value = 592
$ws.Range("J20").Value = 420
$ws.Range("J23").Value = 194
$ws.Range("J25").Value = 102
$ws.Range("J27").Value = 121
$ws.Range("J29").Value = 1144
$ws.Range("J31").Value = 186
$ws.Range("J33").Value = 937
$ws.Range("J36").Value = 281
$ws.Range("J37").Value = 629
$ws.Range("J39").Value = 9
$ws.Range("J42").Value = 840
$ws.Range("H46").Value = 61
$ws.Range("J48").Value = 235
$ws.Range("J52").Value = 521
$ws.Range("J55").Value = 274
$ws.Range("J57").Value = 85
$ws.Range("J60").Value = 125
$ws.Range("J63").Value = 71
$ws.Range("J64").Value = 137
$ws.Range("J65").Value = 525
$ws.Range("J67").Value = 774
$ws.Range("J73").Value = 193
$ws.Range("J76").Value = 295
$ws.Range("J78").Value = 250
$ws.Range("J79").Value = 584
$ws.Range("J83").Value = 418
$ws.Range("J85").Value = 866
$ws.Range("J86").Value = 122
$ws.Range("J89").Value = 265
$ws.Range("J90").Value = 222
$ws.Range("J91").Value = 226
$ws.Range("J93").Value = 94
$ws.Range("J94").Value = 206
$ws.Range("J95").Value = 305
$ws.Range("J97").Value = 168
$ws.Range("J99").Value = 319
$ws.Range("H101").Value = 26013
$ws.Range("J101").Value = 20385

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("J2").Value = 72
$ws.Range("J6").Value = 50
$ws.Range("J7").Value = 186

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("J3").Value = 298
$ws.Range("J4").Value = 60
$ws.Range("J7").Value = 774

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("J3").Value = 76
$ws.Range("J6").Value = 184

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("J2").Value = 342
$ws.Range("J6").Value = 301
$ws.Range("J7").Value = 1144

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("J6").Value = 118
$ws.Range("J7").Value = 235

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("J2").Value = 149
$ws.Range("J3").Value = 176
$ws.Range("J6").Value = 216
$ws.Range("J7").Value = 592

$ws = $wb.Worksheets.Item("River North")
$ws.Range("J6").Value = 163
$ws.Range("J7").Value = 295

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("J2").Value = 39
$ws.Range("J6").Value = 35
$ws.Range("J7").Value = 103

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("J2").Value = 42
$ws.Range("J7").Value = 147

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("J2").Value = 186
$ws.Range("J3").Value = 168
$ws.Range("J6").Value = 431
$ws.Range("J7").Value = 840

$ws = $wb.Worksheets.Item("Boystown")
$ws.Range("J5").Value = 11
$ws.Range("J6").Value = 25

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("J3").Value = 82
$ws.Range("J7").Value = 250

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("J6").Value = 134
$ws.Range("J7").Value = 274

$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("H3").Value = 13
$ws.Range("H7").Value = 61

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("J3").Value = 67
$ws.Range("J7").Value = 194

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("J3").Value = 93
$ws.Range("J7").Value = 226

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("J2").Value = 165
$ws.Range("J3").Value = 206
$ws.Range("J6").Value = 163
$ws.Range("J7").Value = 584

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("J3").Value = 36
$ws.Range("J4").Value = 11
$ws.Range("J7").Value = 137

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("J6").Value = 111
$ws.Range("J7").Value = 420

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("J3").Value = 36
$ws.Range("J7").Value = 173

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("J6").Value = 84
$ws.Range("J7").Value = 281

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("J2").Value = 22
$ws.Range("J7").Value = 94

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("J2").Value = 184
$ws.Range("J3").Value = 180
$ws.Range("J6").Value = 185
$ws.Range("J7").Value = 589

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("J3").Value = 42
$ws.Range("J7").Value = 206

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("J2").Value = 46
$ws.Range("J7").Value = 102

$ws = $wb.Worksheets.Item("Greektown")
$ws.Range("J2").Value = 3
$ws.Range("J6").Value = 9

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("J6").Value = 129
$ws.Range("J7").Value = 318

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("J2").Value = 68
$ws.Range("J7").Value = 193

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("J3").Value = 18
$ws.Range("J7").Value = 168

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("J6").Value = 77
$ws.Range("J7").Value = 265

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("J3").Value = 29
$ws.Range("J7").Value = 121

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("J3").Value = 18
$ws.Range("J7").Value = 122

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("J2").Value = 82
$ws.Range("J6").Value = 63
$ws.Range("J7").Value = 222

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("J3").Value = 24
$ws.Range("J7").Value = 85

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("J3").Value = 36
$ws.Range("J7").Value = 125

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("J2").Value = 226
$ws.Range("J3").Value = 317
$ws.Range("J6").Value = 250
$ws.Range("J7").Value = 866

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("J2").Value = 124
$ws.Range("J3").Value = 159
$ws.Range("J6").Value = 209
$ws.Range("J7").Value = 521
